$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "ofmethod" row (row 2): change the value used for optical-flow method
# from "hs" to "farneback" for the first parameter column (B2), while the
# second parameter column (C2) keeps its previous value ("hs").
$ws.Range("B2").Value = "farneback"
$ws.Range("C2").Value = "hs"
